# Updated cryptos list values (price + 1h volume-change columns) to match the
# latest scrape. Each entry is a cell reference plus its new literal text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '30.870.60' }
    @{ Cell = 'E2'; Value = '  -1.16%  ' }
    @{ Cell = 'D3'; Value = '1.939.14' }
    @{ Cell = 'E3'; Value = '  -1.29%  ' }
    @{ Cell = 'E4'; Value = '  -0.38%  ' }
    @{ Cell = 'D5'; Value = '243.74' }
    @{ Cell = 'E5'; Value = '  -1.09%  ' }
    @{ Cell = 'D6'; Value = '1.000' }
    @{ Cell = 'E6'; Value = '  -0.26%  ' }
    @{ Cell = 'D7'; Value = '0.4915' }
    @{ Cell = 'E7'; Value = '  -0.21%  ' }
    @{ Cell = 'D8'; Value = '0.2945' }
    @{ Cell = 'E8'; Value = '  -2.32%  ' }
    @{ Cell = 'D9'; Value = '0.06890' }
    @{ Cell = 'E9'; Value = '  +0.00%  ' }
    @{ Cell = 'D10'; Value = '19.24' }
    @{ Cell = 'E10'; Value = '  -0.32%  ' }
    @{ Cell = 'D11'; Value = '105.34' }
    @{ Cell = 'E11'; Value = '  -3.48%  ' }
    @{ Cell = 'D12'; Value = '1.935.57' }
    @{ Cell = 'E12'; Value = '  -1.22%  ' }
    @{ Cell = 'D13'; Value = '0.07770' }
    @{ Cell = 'E13'; Value = '  -0.28%  ' }
    @{ Cell = 'D14'; Value = '5.369' }
    @{ Cell = 'E14'; Value = '  -2.00%  ' }
    @{ Cell = 'D15'; Value = '0.7041' }
    @{ Cell = 'E15'; Value = '  -1.50%  ' }
    @{ Cell = 'D16'; Value = '275.38' }
    @{ Cell = 'E16'; Value = '  -3.53%  ' }
    @{ Cell = 'D17'; Value = '30.870.44' }
    @{ Cell = 'E17'; Value = '  -0.77%  ' }
    @{ Cell = 'D18'; Value = '0.000007739' }
    @{ Cell = 'E18'; Value = '  -0.59%  ' }
    @{ Cell = 'D19'; Value = '13.12' }
    @{ Cell = 'E19'; Value = '  -1.54%  ' }
    @{ Cell = 'B20'; Value = 'WrappedliquidstakedEther2.0' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth' }
    @{ Cell = 'D20'; Value = '2.198.21' }
    @{ Cell = 'E20'; Value = '  +0.16%  ' }
    @{ Cell = 'B21'; Value = 'Uniswap' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni' }
    @{ Cell = 'D21'; Value = '5.580' }
    @{ Cell = 'E21'; Value = '  +0.94%  ' }
    @{ Cell = 'B22'; Value = 'Dai' }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai' }
    @{ Cell = 'D22'; Value = '1.000' }
    @{ Cell = 'E22'; Value = '  -0.31%  ' }
    @{ Cell = 'D23'; Value = '1.001' }
    @{ Cell = 'E23'; Value = '  -0.49%  ' }
    @{ Cell = 'D24'; Value = '6.542' }
    @{ Cell = 'E24'; Value = '  -0.66%  ' }
    @{ Cell = 'D25'; Value = '9.856' }
    @{ Cell = 'E25'; Value = '  -0.04%  ' }
    @{ Cell = 'D26'; Value = '166.09' }
    @{ Cell = 'E26'; Value = '  -2.27%  ' }
    @{ Cell = 'D27'; Value = '19.64' }
    @{ Cell = 'E27'; Value = '  -3.08%  ' }
    @{ Cell = 'D28'; Value = '2.155' }
    @{ Cell = 'E28'; Value = '  -3.83%  ' }
    @{ Cell = 'D29'; Value = '0.1044' }
    @{ Cell = 'E29'; Value = '  -0.71%  ' }
    @{ Cell = 'D30'; Value = '1.391' }
    @{ Cell = 'E30'; Value = '  -3.43%  ' }
    @{ Cell = 'D31'; Value = '1.562' }
    @{ Cell = 'E31'; Value = '  -1.54%  ' }
    @{ Cell = 'D32'; Value = '4.568' }
    @{ Cell = 'E32'; Value = '  -1.27%  ' }
    @{ Cell = 'D33'; Value = '4.385' }
    @{ Cell = 'E33'; Value = '  -2.54%  ' }
    @{ Cell = 'D34'; Value = '0.04896' }
    @{ Cell = 'E34'; Value = '  -1.91%  ' }
    @{ Cell = 'D35'; Value = '0.7591' }
    @{ Cell = 'E35'; Value = '  -0.63%  ' }
    @{ Cell = 'D36'; Value = '1.153' }
    @{ Cell = 'E36'; Value = '  -2.79%  ' }
    @{ Cell = 'D37'; Value = '0.9997' }
    @{ Cell = 'E37'; Value = '  -0.31%  ' }
    @{ Cell = 'D38'; Value = '2.736' }
    @{ Cell = 'E38'; Value = '  -0.21%  ' }
    @{ Cell = 'D39'; Value = '0.02009' }
    @{ Cell = 'E39'; Value = '  -2.19%  ' }
    @{ Cell = 'D40'; Value = '2.660' }
    @{ Cell = 'E40'; Value = '  -2.04%  ' }
    @{ Cell = 'B41'; Value = 'FraxShare' }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs' }
    @{ Cell = 'D41'; Value = '6.516' }
    @{ Cell = 'E41'; Value = '  +0.85%  ' }
    @{ Cell = 'B42'; Value = 'Aave' }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave' }
    @{ Cell = 'D42'; Value = '78.10' }
    @{ Cell = 'E42'; Value = '  +6.81%  ' }
    @{ Cell = 'D43'; Value = '2.096' }
    @{ Cell = 'E43'; Value = '  -3.91%  ' }
    @{ Cell = 'D44'; Value = '0.9135' }
    @{ Cell = 'E44'; Value = '  +3.23%  ' }
    @{ Cell = 'D45'; Value = '0.4444' }
    @{ Cell = 'E45'; Value = '  -2.09%  ' }
    @{ Cell = 'D46'; Value = '107.66' }
    @{ Cell = 'E46'; Value = '  -1.81%  ' }
    @{ Cell = 'D47'; Value = '0.9989' }
    @{ Cell = 'E47'; Value = '  -0.53%  ' }
    @{ Cell = 'D48'; Value = '7.678' }
    @{ Cell = 'E48'; Value = '  -6.83%  ' }
    @{ Cell = 'D49'; Value = '998.02' }
    @{ Cell = 'E49'; Value = '  +3.95%  ' }
    @{ Cell = 'D50'; Value = '0.1247' }
    @{ Cell = 'E50'; Value = '  -1.85%  ' }
    @{ Cell = 'D51'; Value = '36.12' }
    @{ Cell = 'E51'; Value = '  +0.91%  ' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Columns D/E hold price/volume text that looks numeric (e.g. "1.000",
    # "30.870.60"); force the cell to Text format first so Excel keeps the
    # exact original string instead of silently coercing it to a Double and
    # dropping significant trailing digits.
    if ($u.Cell -match '^D') {
        $cell.NumberFormat = '@'
    }
    $cell.Value = $u.Value
}
